# Updated cryptos list on Wed Sep  4 08:50:17 UTC 2024 with GitHub Actions
# Applies the per-row Price (D) and Volume(1h) (E) updates plus the
# Mantle/Bittensor row swap (rows 46 & 47) described by the diff.
#
# Note: several new Price values look like plain numbers (e.g. "0.999").
# Excel's COM layer auto-converts such text into a numeric cell when the
# Value is assigned, so for those cells we first force the NumberFormat
# to Text ("@") so the literal string is preserved, matching the
# original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 46 & 47: Bittensor and Mantle swap positions with new values ---
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "255.85"
$ws.Range("E46").Value = "  -6.74%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.571"
$ws.Range("E47").Value = "  -3.15%  "

# --- Remaining Price / Volume(1h) updates ---
$ws.Range("D2").Value = "56.656.36"
$ws.Range("E2").Value = "  -3.50%  "
$ws.Range("D3").Value = "2.374.87"
$ws.Range("E3").Value = "  -4.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.51"
$ws.Range("E5").Value = "  -4.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.97"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").Value = "2.395.60"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0965"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.68"
$ws.Range("E13").Value = "  -9.68%  "
$ws.Range("D14").Value = "2.799.15"
$ws.Range("E14").Value = "  -4.57%  "
$ws.Range("D15").Value = "56.443.81"
$ws.Range("E15").Value = "  -3.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.62"
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "2.388.03"
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.24"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.89"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.06"
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.28"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.69"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "2.496.70"
$ws.Range("E26").Value = "  -4.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.380"
$ws.Range("E27").Value = "  -6.69%  "
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.88"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.68"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("E32").Value = "  -4.77%  "
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("E34").Value = "  -6.00%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.76"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("E39").Value = "  -4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.86"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").Value = "  -4.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.785"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.03"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.89"
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0901"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0488"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.79"
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0208"
$ws.Range("E51").Value = "  -4.52%  "
